$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28,8).Value = 607.5714  # H28: was 721
$ws.Cells.Item(28,9).Value = 607.5714  # I28: was 736.6
$ws.Cells.Item(28,10).Value = 0  # J28: was 701.5
$ws.Cells.Item(28,11).Value = 607.5714  # K28: was 736.6
$ws.Cells.Item(28,12).Value = 0  # L28: was 701.5
$ws.Cells.Item(28,13).Value = -122.5714  # M28: was -251.6
$ws.Cells.Item(28,14).ClearContents()  # N28: was -1671.5

$ws.Cells.Item(80,8).Value = 50000  # H80: was 25199.5
$ws.Cells.Item(80,10).Value = 50000  # J80: was 25199.5
$ws.Cells.Item(80,12).Value = 150000  # L80: was 75598.5
$ws.Cells.Item(80,14).Value = -151996  # N80: was -77594.5

$ws.Cells.Item(83,8).Value = 50000  # H83: was 25199.5
$ws.Cells.Item(83,10).Value = 50000  # J83: was 25199.5
$ws.Cells.Item(83,12).Value = 450000  # L83: was 226795.5
$ws.Cells.Item(83,14).Value = -459984  # N83: was -236779.5

$ws.Cells.Item(130,8).Value = 96657.336  # H130: was 95978
$ws.Cells.Item(130,10).Value = 96657.336  # J130: was 95978
$ws.Cells.Item(130,12).Value = 96657.336  # L130: was 95978
$ws.Cells.Item(130,14).Value = -106697.336  # N130: was -106018

$ws.Cells.Item(132,8).Value = 3087.9  # H132: was 3177.1724
$ws.Cells.Item(132,9).Value = 2426.5417  # I132: was 2510.348
$ws.Cells.Item(132,11).Value = 7279.625100000001  # K132: was 7531.044
$ws.Cells.Item(132,13).Value = -4749.625100000001  # M132: was -5001.044

$ws.Cells.Item(133,8).Value = 99780  # H133: was 99995
$ws.Cells.Item(133,10).Value = 99780  # J133: was 99995
$ws.Cells.Item(133,12).Value = 99780  # L133: was 99995
$ws.Cells.Item(133,14).Value = -109900  # N133: was -110115

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(17,8).Value = 1002  # H17: was 1004
$ws.Cells.Item(17,10).Value = 1002  # J17: was 1004
$ws.Cells.Item(17,12).Value = 1002  # L17: was 1004
$ws.Cells.Item(17,14).Value = -1348  # N17: was -1350

$ws.Cells.Item(61,8).Value = 3583.1667  # H61: was 3625
$ws.Cells.Item(61,9).Value = 3583.1667  # I61: was 3625
$ws.Cells.Item(61,11).Value = 3583.1667  # K61: was 3625
$ws.Cells.Item(61,13).Value = -3371.1667  # M61: was -3413

$ws.Cells.Item(74,8).Value = 2297.1  # H74: was 4399.6
$ws.Cells.Item(74,9).Value = 1567.2858  # I74: was 1666
$ws.Cells.Item(74,10).Value = 4000  # J74: was 8500
$ws.Cells.Item(74,11).Value = 1567.2858  # K74: was 1666
$ws.Cells.Item(74,12).Value = 4000  # L74: was 8500
$ws.Cells.Item(74,13).Value = -693.2858000000001  # M74: was -792
$ws.Cells.Item(74,14).Value = -5748  # N74: was -10248

$ws.Cells.Item(77,8).Value = 2297.1  # H77: was 4399.6
$ws.Cells.Item(77,9).Value = 1567.2858  # I77: was 1666
$ws.Cells.Item(77,10).Value = 4000  # J77: was 8500
$ws.Cells.Item(77,11).Value = 7836.429  # K77: was 8330
$ws.Cells.Item(77,12).Value = 20000  # L77: was 42500
$ws.Cells.Item(77,13).Value = -3468.429  # M77: was -3962
$ws.Cells.Item(77,14).Value = -28736  # N77: was -51236

$ws.Cells.Item(108,8).Value = 75925  # H108: was 75950
$ws.Cells.Item(108,10).Value = 75925  # J108: was 75950
$ws.Cells.Item(108,12).Value = 75925  # L108: was 75950
$ws.Cells.Item(108,14).Value = -83605  # N108: was -83630

$ws.Cells.Item(109,8).Value = 0  # H109: was 25000
$ws.Cells.Item(109,9).Value = 0  # I109: was 25000
$ws.Cells.Item(109,11).Value = 0  # K109: was 25000
$ws.Cells.Item(109,13).ClearContents()  # M109: was -23613

$ws.Cells.Item(111,8).Value = 99999  # H111: was 99995
$ws.Cells.Item(111,10).Value = 99999  # J111: was 99995
$ws.Cells.Item(111,12).Value = 99999  # L111: was 99995
$ws.Cells.Item(111,14).Value = -108179  # N111: was -108175

$ws.Cells.Item(125,8).Value = 39500  # H125: was 0
$ws.Cells.Item(125,10).Value = 39500  # J125: was 0
$ws.Cells.Item(125,12).Value = 39500  # L125: was 0
$ws.Cells.Item(125,14).Value = -49340  # N125: was None

$ws.Cells.Item(135,8).Value = 70000  # H135: was 85000
$ws.Cells.Item(135,10).Value = 70000  # J135: was 85000
$ws.Cells.Item(135,12).Value = 70000  # L135: was 85000
$ws.Cells.Item(135,14).Value = -80140  # N135: was -95140

$ws.Cells.Item(136,8).Value = 3583.1667  # H136: was 3625
$ws.Cells.Item(136,9).Value = 3583.1667  # I136: was 3625
$ws.Cells.Item(136,11).Value = 10749.5001  # K136: was 10875
$ws.Cells.Item(136,13).Value = -8199.500100000001  # M136: was -8325

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80,8).Value = 582.8333  # H80: was 520.4
$ws.Cells.Item(80,10).Value = 833.6667  # J80: was 803
$ws.Cells.Item(80,12).Value = 833.6667  # L80: was 803
$ws.Cells.Item(80,14).Value = -2829.6667  # N80: was -2799

$ws.Cells.Item(83,8).Value = 582.8333  # H83: was 520.4
$ws.Cells.Item(83,10).Value = 833.6667  # J83: was 803
$ws.Cells.Item(83,12).Value = 4168.3335  # L83: was 4015
$ws.Cells.Item(83,14).Value = -14152.3335  # N83: was -13999

$ws.Cells.Item(115,8).Value = 0  # H115: was 90000
$ws.Cells.Item(115,9).Value = 0  # I115: was 90000
$ws.Cells.Item(115,11).Value = 0  # K115: was 90000
$ws.Cells.Item(115,13).ClearContents()  # M115: was -88433

$ws.Cells.Item(116,8).Value = 85000  # H116: was 92497.5
$ws.Cells.Item(116,10).Value = 85000  # J116: was 92497.5
$ws.Cells.Item(116,12).Value = 85000  # L116: was 92497.5
$ws.Cells.Item(116,14).Value = -94178  # N116: was -101675.5

$ws.Cells.Item(137,8).Value = 0  # H137: was 99995
$ws.Cells.Item(137,10).Value = 0  # J137: was 99995
$ws.Cells.Item(137,12).Value = 0  # L137: was 99995
$ws.Cells.Item(137,14).ClearContents()  # N137: was -110195

$ws.Cells.Item(140,8).Value = 0  # H140: was 99995
$ws.Cells.Item(140,10).Value = 0  # J140: was 99995
$ws.Cells.Item(140,12).Value = 0  # L140: was 99995
$ws.Cells.Item(140,14).ClearContents()  # N140: was -110355

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31,8).Value = 4159.375  # H31: was 3997.7646
$ws.Cells.Item(31,9).Value = 1997  # I31: was 1913.4286
$ws.Cells.Item(31,11).Value = 1997  # K31: was 1913.4286
$ws.Cells.Item(31,13).Value = -1702  # M31: was -1618.4286

$ws.Cells.Item(34,8).Value = 4159.375  # H34: was 3997.7646
$ws.Cells.Item(34,9).Value = 1997  # I34: was 1913.4286
$ws.Cells.Item(34,11).Value = 1997  # K34: was 1913.4286
$ws.Cells.Item(34,13).Value = -1795  # M34: was -1711.4286

$ws.Cells.Item(53,8).Value = 0  # H53: was 99995
$ws.Cells.Item(53,10).Value = 0  # J53: was 99995
$ws.Cells.Item(53,12).Value = 0  # L53: was 99995
$ws.Cells.Item(53,14).ClearContents()  # N53: was -101209

$ws.Cells.Item(87,8).Value = 33332.668  # H87: was 32500
$ws.Cells.Item(87,9).Value = 0  # I87: was 25000
$ws.Cells.Item(87,10).Value = 33332.668  # J87: was 40000
$ws.Cells.Item(87,11).Value = 0  # K87: was 25000
$ws.Cells.Item(87,12).Value = 33332.668  # L87: was 40000
$ws.Cells.Item(87,13).ClearContents()  # M87: was -23814
$ws.Cells.Item(87,14).Value = -35704.668  # N87: was -42372

$ws.Cells.Item(90,8).Value = 33332.668  # H90: was 32500
$ws.Cells.Item(90,9).Value = 0  # I90: was 25000
$ws.Cells.Item(90,10).Value = 33332.668  # J90: was 40000
$ws.Cells.Item(90,11).Value = 0  # K90: was 75000
$ws.Cells.Item(90,12).Value = 99998.00399999999  # L90: was 120000
$ws.Cells.Item(90,13).ClearContents()  # M90: was -69072
$ws.Cells.Item(90,14).Value = -111854.004  # N90: was -131856

$ws.Cells.Item(111,8).Value = 0  # H111: was 99995
$ws.Cells.Item(111,10).Value = 0  # J111: was 99995
$ws.Cells.Item(111,12).Value = 0  # L111: was 99995
$ws.Cells.Item(111,14).ClearContents()  # N111: was -108175

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(31,8).Value = 750  # H31: was 0
$ws.Cells.Item(31,9).Value = 700  # I31: was 0
$ws.Cells.Item(31,10).Value = 800  # J31: was 0
$ws.Cells.Item(31,11).Value = 2100  # K31: was 0
$ws.Cells.Item(31,12).Value = 2400  # L31: was 0
$ws.Cells.Item(31,13).Value = -1812  # M31: was None
$ws.Cells.Item(31,14).Value = -2976  # N31: was None

$ws.Cells.Item(68,8).Value = 1349.5  # H68: was 1533
$ws.Cells.Item(68,9).Value = 799.5  # I68: was 800
$ws.Cells.Item(68,11).Value = 2398.5  # K68: was 2400
$ws.Cells.Item(68,13).Value = -1587.5  # M68: was -1589

$ws.Cells.Item(71,8).Value = 1349.5  # H71: was 1533
$ws.Cells.Item(71,9).Value = 799.5  # I71: was 800
$ws.Cells.Item(71,11).Value = 7195.5  # K71: was 7200
$ws.Cells.Item(71,13).Value = -3139.5  # M71: was -3144

$ws.Cells.Item(114,8).Value = 729.5  # H114: was 528
$ws.Cells.Item(114,10).Value = 931  # J114: was 0
$ws.Cells.Item(114,12).Value = 2793  # L114: was 0
$ws.Cells.Item(114,14).Value = -9301  # N114: was None

$ws.Cells.Item(137,8).Value = 4600  # H137: was 3959.8
$ws.Cells.Item(137,9).Value = 4600  # I137: was 3959.8
$ws.Cells.Item(137,11).Value = 13800  # K137: was 11879.4
$ws.Cells.Item(137,13).Value = -8700  # M137: was -6779.400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(12,8).Value = 2260.6  # H12: was 2750
$ws.Cells.Item(12,9).Value = 2260.6  # I12: was 2750
$ws.Cells.Item(12,11).Value = 2260.6  # K12: was 2750
$ws.Cells.Item(12,13).Value = -2120.6  # M12: was -2610

$ws.Cells.Item(48,8).Value = 13333.333  # H48: was 27500
$ws.Cells.Item(48,10).Value = 13333.333  # J48: was 27500
$ws.Cells.Item(48,12).Value = 13333.333  # L48: was 27500
$ws.Cells.Item(48,14).Value = -14303.333  # N48: was -28470

$ws.Cells.Item(80,8).Value = 3248.75  # H80: was 3416.6667
$ws.Cells.Item(80,9).Value = 2831.6667  # I80: was 2875
$ws.Cells.Item(80,11).Value = 2831.6667  # K80: was 2875
$ws.Cells.Item(80,13).Value = -1833.6667  # M80: was -1877

$ws.Cells.Item(83,8).Value = 3248.75  # H83: was 3416.6667
$ws.Cells.Item(83,9).Value = 2831.6667  # I83: was 2875
$ws.Cells.Item(83,11).Value = 14158.3335  # K83: was 14375
$ws.Cells.Item(83,13).Value = -9166.333500000001  # M83: was -9383

$ws.Cells.Item(124,8).Value = 0  # H124: was 99995
$ws.Cells.Item(124,10).Value = 0  # J124: was 99995
$ws.Cells.Item(124,12).Value = 0  # L124: was 99995
$ws.Cells.Item(124,14).ClearContents()  # N124: was -109815

$ws.Cells.Item(132,8).Value = 7544.909  # H132: was 8221.556
$ws.Cells.Item(132,9).Value = 6936.875  # I132: was 7749.1665
$ws.Cells.Item(132,11).Value = 20810.625  # K132: was 23247.4995
$ws.Cells.Item(132,13).Value = -18280.625  # M132: was -20717.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22,8).Value = 1656.619  # H22: was 1599.2778
$ws.Cells.Item(22,9).Value = 1599.375  # I22: was 1506.7693
$ws.Cells.Item(22,11).Value = 1599.375  # K22: was 1506.7693
$ws.Cells.Item(22,13).Value = -1304.375  # M22: was -1211.7693

$ws.Cells.Item(25,8).Value = 10000  # H25: was 0
$ws.Cells.Item(25,10).Value = 10000  # J25: was 0
$ws.Cells.Item(25,12).Value = 10000  # L25: was 0
$ws.Cells.Item(25,14).Value = -10460  # N25: was None

$ws.Cells.Item(27,8).Value = 1656.619  # H27: was 1599.2778
$ws.Cells.Item(27,9).Value = 1599.375  # I27: was 1506.7693
$ws.Cells.Item(27,11).Value = 1599.375  # K27: was 1506.7693
$ws.Cells.Item(27,13).Value = -1492.375  # M27: was -1399.7693

$ws.Cells.Item(46,8).Value = 4714.7646  # H46: was 4615.091
$ws.Cells.Item(46,10).Value = 4855.2905  # J46: was 4750.3335
$ws.Cells.Item(46,12).Value = 4855.2905  # L46: was 4750.3335
$ws.Cells.Item(46,14).Value = -5231.2905  # N46: was -5126.3335

$ws.Cells.Item(68,8).Value = 2771  # H68: was 2831.1333
$ws.Cells.Item(68,9).Value = 2433.2144  # I68: was 2497.6428
$ws.Cells.Item(68,11).Value = 2433.2144  # K68: was 2497.6428
$ws.Cells.Item(68,13).Value = -1684.2144  # M68: was -1748.6428

$ws.Cells.Item(71,8).Value = 2771  # H71: was 2831.1333
$ws.Cells.Item(71,9).Value = 2433.2144  # I71: was 2497.6428
$ws.Cells.Item(71,11).Value = 12166.072  # K71: was 12488.214
$ws.Cells.Item(71,13).Value = -8422.072  # M71: was -8744.214

$ws.Cells.Item(124,8).Value = 99990  # H124: was 0
$ws.Cells.Item(124,10).Value = 99990  # J124: was 0
$ws.Cells.Item(124,12).Value = 99990  # L124: was 0
$ws.Cells.Item(124,14).Value = -109810  # N124: was None

$ws.Cells.Item(138,8).Value = 99429  # H138: was 0
$ws.Cells.Item(138,10).Value = 99429  # J138: was 0
$ws.Cells.Item(138,12).Value = 99429  # L138: was 0
$ws.Cells.Item(138,14).Value = -109709  # N138: was None

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(48,8).Value = 20000  # H48: was 0
$ws.Cells.Item(48,10).Value = 20000  # J48: was 0
$ws.Cells.Item(48,12).Value = 20000  # L48: was 0
$ws.Cells.Item(48,14).Value = -21138  # N48: was None

$ws.Cells.Item(132,8).Value = 982.3333  # H132: was 984
$ws.Cells.Item(132,9).Value = 982.3333  # I132: was 984
$ws.Cells.Item(132,11).Value = 2946.9999  # K132: was 2952
$ws.Cells.Item(132,13).Value = -416.9998999999998  # M132: was -422

$ws.Cells.Item(136,8).Value = 4002.8333  # H136: was 4464.154
$ws.Cells.Item(136,10).Value = 0  # J136: was 10000
$ws.Cells.Item(136,12).Value = 0  # L136: was 30000
$ws.Cells.Item(136,14).ClearContents()  # N136: was -35100
